$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new column D (shifts the old "Who writes it" column to E) ---
$ws.Columns("D").Insert()

# New column D header + page-count values
$ws.Range("D3").Value = "Maximun number of pages (whole chapter)"
$ws.Range("D4").Value = 1
$ws.Range("D5").Value = 4
$ws.Range("D9").Value = 2
$ws.Range("D11").Value = 1

# Column D width (matches width="37.5" after Excel's internal padding offset)
$ws.Range("D1").ColumnWidth = 36.666666666666664

# --- 2. Row heights: thick-separator rows get a slightly taller row (17pt) ---
$ws.Rows("2:2").RowHeight = 17
$ws.Rows("3:3").RowHeight = 17
$ws.Rows("4:4").RowHeight = 17
$ws.Rows("8:8").RowHeight = 17
$ws.Rows("10:10").RowHeight = 17
$ws.Rows("11:11").RowHeight = 17

# --- 3. Borders: draw a medium box around the whole table (A3:E11) with   ---
#        internal separators after the header (row3/4) and around each    ---
#        chapter block (rows 5-8, 9-10).
$ws.Range("A3:E3").Borders.Item(8).Weight = -4138   # top
$ws.Range("A3:E3").Borders.Item(9).Weight = -4138   # bottom
$ws.Range("A4:E4").Borders.Item(8).Weight = -4138   # top
$ws.Range("A4:E4").Borders.Item(9).Weight = -4138   # bottom
$ws.Range("A5:E5").Borders.Item(8).Weight = -4138   # top
$ws.Range("A8:E8").Borders.Item(9).Weight = -4138   # bottom
$ws.Range("A9:E9").Borders.Item(8).Weight = -4138   # top
$ws.Range("A10:E10").Borders.Item(9).Weight = -4138 # bottom
$ws.Range("A11:E11").Borders.Item(8).Weight = -4138 # top
$ws.Range("A11:E11").Borders.Item(9).Weight = -4138 # bottom
$ws.Range("A3:A11").Borders.Item(7).Weight = -4138  # left outline
$ws.Range("E3:E11").Borders.Item(10).Weight = -4138 # right outline

# --- 4. Selection / dimension bookkeeping ---
$ws.Range("D12").Select()
